$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("workflow_repository_tools")

# Insert a new row at position 2 (existing rows 2..20 shift down to 3..21),
# then fill in the new "NGSCheckMate" tool entry that documents the
# NBL (neuroblastoma) metastases snippet's associated workflow.
$ws.Rows.Item(2).Insert()

$ws.Cells.Item(2, 1).Value = "NGSCheckMate"

# "1.3" parses as a number via .Value, which would store it as a numeric
# cell (and pull in a new "@" text style not present in the target file).
# Routing it through a literal-string formula, then collapsing that
# formula to its cached value via copy/paste-values, keeps it as a plain
# shared-string cell with no style change - matching how the version
# column stores every other entry in this sheet (e.g. "11.6", "0.9.3").
$verCell = $ws.Cells.Item(2, 2)
$verCell.Formula = "=""1.3"""
$verCell.Copy()
$verCell.PasteSpecial(-4163)

$ws.Cells.Item(2, 3).Value = "https://github.com/d3b-center/OpenPBTA-workflows/blob/master/cwl/bcf_call.cwl"
